$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Thesaurus")

# "unité expérimentale" -> "parcelle unitaire" (renaming the thesaurus entry)
$ws.Range("B12").Value = "parcelle unitaire"
$ws.Range("C12").Value = "parcelle élémentaire,unité expérimentale"
$ws.Range("E12").Value = "La parcelle unitaire (ou élémentaire) est la plus petite unité expérimentale qui reçoit un traitement. Dans les essais au champ, elle est constituée d'une certaine étendue de terrain et d'un certain nombre de plants. "

# Remove the now-obsolete INRAE thesaurus link cell (and its hyperlink)
$ws.Range("J12").Hyperlinks.Delete()
$ws.Range("J12").ClearContents()
